# Auto-generated Excel COM-interop script to apply the forecast data bugfix.
# The edit:
#   - Shifts all forecast rows up by one (the original row 2 "2007" stub row is removed)
#   - Recomputes the forecasted C (y_1) and E (y_1_forecast) columns with corrected values
#   - Removes the now-stale last row (53), shrinking the used range to A1:E52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(39583, 2008, $null, 2009, $null),
    @(39765, 2008, $null, 2009, $null),
    @(39948, 2009, $null, 2010, $null),
    @(40130, 2009, -1.317619545389281, 2010, $null),
    @(40310, 2010, $null, 2011, $null),
    @(40494, 2010, 2.771597318554297, 2011, $null),
    @(40676, 2011, $null, 2012, $null),
    @(40862, 2011, 1.799362536952542, 2012, $null),
    @(41044, 2012, $null, 2013, $null),
    @(41228, 2012, 2.123182427147152, 2013, 4.792854588620821),
    @(41409, 2013, 4.880442637054072, 2014, 1.641301872652501),
    @(41592, 2013, 7.317297369612819, 2014, 7.239454936865775),
    @(41774, 2014, 5.941867202078877, 2015, 2.672847571394987),
    @(41957, 2014, 4.260319658857736, 2015, 0.4361429468412448),
    @(42137, 2015, 0.292749233164491, 2016, 2.181874035977249),
    @(42321, 2015, -0.05262415810141086, 2016, 1.392195163617171),
    @(42503, 2016, 2.032207428223742, 2017, 2.13692496326825),
    @(42689, 2016, 1.459778471779982, 2017, 1.493220091771108),
    @(42867, 2017, 2.775332754349846, 2018, 2.200426660963761),
    @(43053, 2017, 3.002208343813528, 2018, 2.622364272988187),
    @(43145, 2018, 3.2651197821016, 2019, 2.271936475508851),
    @(43235, 2018, 3.565025829754953, 2019, 2.446228176258058),
    @(43326, 2018, 4.01493878081518, 2019, 3.091110147865495),
    @(43418, 2018, 4.020433260014977, 2019, 3.113086948791377),
    @(43510, 2019, 3.283136334808323, 2020, 2.403408536719187),
    @(43600, 2019, 3.444206290325491, 2020, 2.496958452261078),
    @(43691, 2019, 3.479628752085517, 2020, 2.54748094003614),
    @(43783, 2019, 3.53224976671227, 2020, 2.757421718286168),
    @(43875, 2020, 2.31260691849986, 2021, 2.363182008239928),
    @(43966, 2020, 2.667234932970275, 2021, 2.567662999186382),
    @(44068, 2020, -0.985458715495402, 2021, -1.999977666418695),
    @(44159, 2020, -0.985458715495402, 2021, -1.373617952268746),
    @(44251, 2021, -4.853362183897836, 2022, 0.6599234717970859),
    @(44341, 2021, -4.511102905979703, 2022, 0.9582724917052587),
    @(44432, 2021, -4.365687260408224, 2022, 1.582150300418306),
    @(44525, 2021, -4.365687260408224, 2022, 1.324283050325015),
    @(44617, 2022, 2.700663803921799, 2023, 2.813733698830267),
    @(44706, 2022, 1.386772772629241, 2023, 0.8813242377093244),
    @(44798, 2022, 1.27347919322387, 2023, 0.4287685802702779),
    @(44890, 2022, 1.27347919322387, 2023, -0.5025420863900898),
    @(44981, 2023, -0.185315122156382, 2024, 1.715279842342743),
    @(45071, 2023, -0.9537175292835154, 2024, 0.9049225073274991),
    @(45163, 2023, -1.339436245206127, 2024, -0.6598239038080322),
    @(45254, 2023, -1.339436245206127, 2024, -2.371854438773213),
    @(45345, 2024, -2.5174493871855, 2025, 1.089612584330668),
    @(45436, 2024, -3.303819519576723, 2025, -0.3655818470008065),
    @(45534, 2024, -3.451527003230626, 2025, -1.070961900287937),
    @(45618, 2024, -3.451527003230626, 2025, -1.436963918858969),
    @(45713, 2025, -1.652703512303566, 2026, 0.4469810487905734),
    @(45800, 2025, -1.704805397136089, 2026, 0.1460701281005727),
    @(45891, 2025, -1.75044229618867, 2026, -0.1564272439962933)
)

$rowCount = $values.Count
$colCount = 5

$data = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $data[$i,$j] = $values[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 5))
$rng.Value = $data

# Remove the now-obsolete trailing row (previously row 53) so the sheet ends at row 52.
$ws.Rows($endRow + 1).Delete()

Write-Host "Applied forecast series update; rows now span 2..$endRow"
